$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 216.6
$ws.Range("J9").Value = 500
$ws.Range("L9").Value = 500
$ws.Range("N9").Value = -838
$ws.Range("H125").Value = 2995
$ws.Range("I125").Value = 2995
$ws.Range("K125").Value = 26955
$ws.Range("M125").Value = -24495
$ws.Range("H135").Value = 1507
$ws.Range("I135").Value = 760.75
$ws.Range("J135").Value = 2999.5
$ws.Range("K135").Value = 6846.75
$ws.Range("L135").Value = 26995.5
$ws.Range("M135").Value = -4311.75
$ws.Range("N135").Value = -32065.5
$ws.Range("H141").Value = 4230.4165
$ws.Range("I141").Value = 3926.5
$ws.Range("J141").Value = 5750
$ws.Range("K141").Value = 11779.5
$ws.Range("L141").Value = 17250
$ws.Range("M141").Value = -6599.5
$ws.Range("N141").Value = -27610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 706.8421
$ws.Range("I2").Value = 619.41174
$ws.Range("K2").Value = 619.41174
$ws.Range("M2").Value = -506.41174
$ws.Range("H32").Value = 19368.553
$ws.Range("I32").Value = 9930.25
$ws.Range("K32").Value = 9930.25
$ws.Range("M32").Value = -9643.25
$ws.Range("H61").Value = 1407.8524
$ws.Range("I61").Value = 1403.2587
$ws.Range("K61").Value = 1403.2587
$ws.Range("M61").Value = -1191.2587
$ws.Range("H74").Value = 4415.3125
$ws.Range("I74").Value = 1089.1111
$ws.Range("J74").Value = 8691.857
$ws.Range("K74").Value = 1089.1111
$ws.Range("L74").Value = 8691.857
$ws.Range("M74").Value = -215.1111000000001
$ws.Range("N74").Value = -10439.857
$ws.Range("H77").Value = 4415.3125
$ws.Range("I77").Value = 1089.1111
$ws.Range("J77").Value = 8691.857
$ws.Range("K77").Value = 5445.5555
$ws.Range("L77").Value = 43459.285
$ws.Range("M77").Value = -1077.5555
$ws.Range("N77").Value = -52195.285
$ws.Range("H116").Value = 706.8421
$ws.Range("I116").Value = 619.41174
$ws.Range("K116").Value = 619.41174
$ws.Range("M116").Value = 1674.58826
$ws.Range("H132").Value = 1013.76
$ws.Range("I132").Value = 1014.375
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 3043.125
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -513.125
$ws.Range("N132").Value = -8057
$ws.Range("H136").Value = 1407.8524
$ws.Range("I136").Value = 1403.2587
$ws.Range("K136").Value = 4209.7761
$ws.Range("M136").Value = -1659.7761

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 706.8421
$ws.Range("I3").Value = 619.41174
$ws.Range("K3").Value = 619.41174
$ws.Range("M3").Value = -505.41174
$ws.Range("H134").Value = 3014.7646
$ws.Range("I134").Value = 3049.2258
$ws.Range("K134").Value = 9147.6774
$ws.Range("M134").Value = -6612.6774

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5315.6924
$ws.Range("I31").Value = 2793.1667
$ws.Range("J31").Value = 7477.857
$ws.Range("K31").Value = 2793.1667
$ws.Range("L31").Value = 7477.857
$ws.Range("M31").Value = -2498.1667
$ws.Range("N31").Value = -8067.857
$ws.Range("H34").Value = 5315.6924
$ws.Range("I34").Value = 2793.1667
$ws.Range("J34").Value = 7477.857
$ws.Range("K34").Value = 2793.1667
$ws.Range("L34").Value = 7477.857
$ws.Range("M34").Value = -2591.1667
$ws.Range("N34").Value = -7881.857
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15368
$ws.Range("H86").Value = 11588.875
$ws.Range("J86").Value = 17251.25
$ws.Range("L86").Value = 17251.25
$ws.Range("N86").Value = -19497.25
$ws.Range("H89").Value = 11588.875
$ws.Range("J89").Value = 17251.25
$ws.Range("L89").Value = 86256.25
$ws.Range("N89").Value = -97488.25
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490
$ws.Range("H109").Value = 23972.5
$ws.Range("J109").Value = 23972.5
$ws.Range("L109").Value = 23972.5
$ws.Range("N109").Value = -26052.5
$ws.Range("H132").Value = 2730.5715
$ws.Range("I132").Value = 2171.1904
$ws.Range("K132").Value = 6513.5712
$ws.Range("M132").Value = -3983.5712
$ws.Range("H134").Value = 5410.8
$ws.Range("I134").Value = 3500
$ws.Range("J134").Value = 5888.5
$ws.Range("K134").Value = 10500
$ws.Range("L134").Value = 17665.5
$ws.Range("M134").Value = -7965
$ws.Range("N134").Value = -22735.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4569.3887
$ws.Range("I34").Value = 7390.1113
$ws.Range("J34").Value = 1748.6666
$ws.Range("K34").Value = 22170.3339
$ws.Range("L34").Value = 5245.9998
$ws.Range("M34").Value = -22086.3339
$ws.Range("N34").Value = -5413.9998
$ws.Range("H55").Value = 102705
$ws.Range("J55").Value = 4408.3335
$ws.Range("L55").Value = 13225.0005
$ws.Range("N55").Value = -13579.0005
$ws.Range("H63").Value = 1250
$ws.Range("I63").Value = 1250
$ws.Range("K63").Value = 3750
$ws.Range("M63").Value = -3001
$ws.Range("H66").Value = 1250
$ws.Range("I66").Value = 1250
$ws.Range("K66").Value = 11250
$ws.Range("M66").Value = -7506
$ws.Range("H109").Value = 846.3333
$ws.Range("I109").Value = 415.8
$ws.Range("K109").Value = 1247.4
$ws.Range("M109").Value = -207.4000000000001
$ws.Range("H117").Value = 1567.375
$ws.Range("J117").Value = 1789.8334
$ws.Range("L117").Value = 5369.5002
$ws.Range("N117").Value = -12253.5002
$ws.Range("H127").Value = 1750
$ws.Range("J127").Value = 1750
$ws.Range("L127").Value = 5250
$ws.Range("N127").Value = -15170
$ws.Range("H139").Value = 2912.2666
$ws.Range("I139").Value = 2912.2666
$ws.Range("K139").Value = 8736.799800000001
$ws.Range("M139").Value = -3596.799800000001
$ws.Range("H140").Value = 2746.7693
$ws.Range("I140").Value = 2746.7693
$ws.Range("K140").Value = 8240.3079
$ws.Range("M140").Value = -3060.3079

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 15500
$ws.Range("I53").Value = 1000
$ws.Range("K53").Value = 1000
$ws.Range("M53").Value = -369
$ws.Range("H80").Value = 2815.5454
$ws.Range("I80").Value = 3346.1667
$ws.Range("J80").Value = 2178.8
$ws.Range("K80").Value = 3346.1667
$ws.Range("L80").Value = 2178.8
$ws.Range("M80").Value = -2348.1667
$ws.Range("N80").Value = -4174.8
$ws.Range("H83").Value = 2815.5454
$ws.Range("I83").Value = 3346.1667
$ws.Range("J83").Value = 2178.8
$ws.Range("K83").Value = 16730.8335
$ws.Range("L83").Value = 10894
$ws.Range("M83").Value = -11738.8335
$ws.Range("N83").Value = -20878
$ws.Range("H132").Value = 3320.6667
$ws.Range("I132").Value = 2885.158
$ws.Range("K132").Value = 8655.474
$ws.Range("M132").Value = -6125.474
$ws.Range("H136").Value = 25703.455
$ws.Range("J136").Value = 25703.455
$ws.Range("L136").Value = 77110.36500000001
$ws.Range("N136").Value = -82210.36500000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1948.5
$ws.Range("I7").Value = 1998.2858
$ws.Range("K7").Value = 1998.2858
$ws.Range("M7").Value = -1886.2858
$ws.Range("H33").Value = 238336000
$ws.Range("I33").Value = 238336000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 238336000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -238335710
$ws.Range("N33").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 1948.5
$ws.Range("I126").Value = 1998.2858
$ws.Range("K126").Value = 5994.857400000001
$ws.Range("M126").Value = -3524.857400000001
$ws.Range("H132").Value = 4842
$ws.Range("I132").Value = 4412.852
$ws.Range("K132").Value = 13238.556
$ws.Range("M132").Value = -10708.556
$ws.Range("H136").Value = 3218.4546
$ws.Range("I136").Value = 3044.889
$ws.Range("K136").Value = 9134.667000000001
$ws.Range("M136").Value = -6584.667000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 12333.333
$ws.Range("J32").Value = 12333.333
$ws.Range("L32").Value = 12333.333
$ws.Range("N32").Value = -12967.333
$ws.Range("H132").Value = 1455.7858
$ws.Range("I132").Value = 1529.875
$ws.Range("J132").Value = 1357
$ws.Range("K132").Value = 1357
$ws.Range("L132").Value = 4071
$ws.Range("M132").Value = -2059.625
$ws.Range("N132").Value = -9131
